# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the hyperlink that sits on A40 (keeps the cell's text value intact).
$ws.Range("A40").Hyperlinks.Delete()

# 2. Insert a new row above the old row 39 ("US Census Bureau, 2007").
#    This pushes everything from row 39 downward by one row, so the old
#    row 44/45 pair ("USCB" / long citation) becomes row 45/46, and the
#    dimension grows from A1:E45 to A1:E46.
$ws.Range("A39").EntireRow.Insert()

# 3. Now fix up the values in the affected block so the final layout is:
#      A38 Source:
#      A39 (blank)
#      A40 US Census Bureau, 2007
#      A41 (blank)
#      A42 http://www.census.gov/econ/islandareas/2007/historical_data_2007.html
#      A45 USCB
#      A46 USCB
$ws.Range("A39").Value = ""
$ws.Range("A40").Value = "US Census Bureau, 2007"
$ws.Range("A41").Value = ""
$ws.Range("A42").Value = "http://www.census.gov/econ/islandareas/2007/historical_data_2007.html"
$ws.Range("A45").Value = "USCB"
$ws.Range("A46").Value = "USCB"

# 4. Re-apply the "source" (italic) style to the cells that should carry it,
#    and the "title" (bold) style to the USCB label row, matching the
#    original cellStyles used in this sheet. A40 previously used the
#    "HyperLink" style (underlined/blue) - switch it back to italic/black
#    now that the hyperlink itself has been removed.
$ws.Range("A38").Font.Italic = $true
$ws.Range("A39").Font.Italic = $true
$ws.Range("A40").Font.Italic = $true
$ws.Range("A40").Font.Underline = $false
$ws.Range("A41").Font.Italic = $true
$ws.Range("A42").Font.Italic = $true
$ws.Range("A45").Font.Bold = $true
$ws.Range("A46").Font.Italic = $true
